$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.765.83'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.943.75'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.42%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.11'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.18'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.67%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.943.01'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.45%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.505'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.05'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.57%  '

$ws.Range("E11").Value = '  +5.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000233'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +4.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.49'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.22%  '

$ws.Range("E15").Value = '  -1.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.430.53'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.744.13'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.68'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.939.95'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '437.95'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.17%  '

$ws.Range("E21").Value = '  -1.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.664'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.20'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.70%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.69'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.88%  '

$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.12'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.44%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.77'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.39%  '

$ws.Range("E28").Value = '  -0.02%  '

$ws.Range("E29").Value = '  +0.91%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.31'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +6.17%  '

$ws.Range("E31").Value = '  +0.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0₃0979'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +11.94%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.32'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.48%  '

$ws.Range("E34").Value = '  -0.54%  '

$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.989'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.60'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.25%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.01'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.58'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.39%  '

$ws.Range("E40").Value = '  +1.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.44'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.64%  '

$ws.Range("E42").Value = '  -3.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.279'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.50'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.94%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.688.95'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '134.91'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '356.76'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.20%  '

$ws.Range("E50").Value = '  -0.60%  '

$ws.Range("E51").Value = '  -4.03%  '
